$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 24,14
$data[0,0] = 0.4282382573302925
$data[0,1] = 0.0441531750018811
$data[0,2] = 0.6556695097541194
$data[0,3] = 0.2675544873670006
$data[0,4] = 0
$data[0,5] = 0.8695349742738046
$data[0,6] = 0.9352458498722882
$data[0,7] = 0
$data[0,8] = 0.1393602964036802
$data[0,9] = 0.3954657513750135
$data[0,10] = 0
$data[0,11] = 0
$data[0,12] = 0
$data[0,13] = 3.639368077770797
$data[1,0] = 0.3881058685289815
$data[1,1] = 0.03886133666514979
$data[1,2] = 0.6451922757689488
$data[1,3] = 0.2624194397106621
$data[1,4] = 0
$data[1,5] = 0.8738792508109654
$data[1,6] = 0.9417142787176545
$data[1,7] = 0
$data[1,8] = 0.1359607556885578
$data[1,9] = 0.3526853071757898
$data[1,10] = 0
$data[1,11] = 0
$data[1,12] = 0
$data[1,13] = 3.661814556499749
$data[2,0] = 0.3635260368729973
$data[2,1] = 0.03559633351575542
$data[2,2] = 0.6390897997721936
$data[2,3] = 0.2594075812402252
$data[2,4] = 0
$data[2,5] = 0.8770905960700546
$data[2,6] = 0.9460883387850814
$data[2,7] = 0
$data[2,8] = 0.1339512140312564
$data[2,9] = 0.326422754597246
$data[2,10] = 0
$data[2,11] = 0
$data[2,12] = 0
$data[2,13] = 3.677581700991865
$data[3,0] = 0.3535256216705989
$data[3,1] = 0.03426192844905529
$data[3,2] = 0.6366862539041165
$data[3,3] = 0.2582157394604039
$data[3,4] = 0
$data[3,5] = 0.8785359082891162
$data[3,6] = 0.9479720378956955
$data[3,7] = 0
$data[3,8] = 0.1331518760748978
$data[3,9] = 0.3157223678636569
$data[3,10] = 0
$data[3,11] = 0
$data[3,12] = 0
$data[3,13] = 3.684505900458532
$data[4,0] = 0.3518660475656645
$data[4,1] = 0.03404011898477677
$data[4,2] = 0.6362921797546903
$data[4,3] = 0.2580199812145025
$data[4,4] = 0
$data[4,5] = 0.8787841520034476
$data[4,6] = 0.9482909402539264
$data[4,7] = 0
$data[4,8] = 0.1330203285630631
$data[4,9] = 0.3139457044435972
$data[4,10] = 0
$data[4,11] = 0
$data[4,12] = 0
$data[4,13] = 3.685685789721362
$data[5,0] = 0.3633911019653908
$data[5,1] = 0.03557835290270361
$data[5,2] = 0.63905704740651
$data[5,3] = 0.2593913637939664
$data[5,4] = 0
$data[5,5] = 0.8771095348677491
$data[5,6] = 0.9461133330620797
$data[5,7] = 0
$data[5,8] = 0.1339403546473008
$data[5,9] = 0.3262784372778071
$data[5,10] = 0
$data[5,11] = 0
$data[5,12] = 0
$data[5,13] = 3.677673063120224
$data[6,0] = 0.4143881504529645
$data[6,1] = 0.04233187414858719
$data[6,2] = 0.6519884211414251
$data[6,3] = 0.2657546625068861
$data[6,4] = 0
$data[6,5] = 0.8709199148874518
$data[6,6] = 0.9373926674729347
$data[6,7] = 0
$data[6,8] = 0.1381719913948203
$data[6,9] = 0.3807144048705595
$data[6,10] = 0
$data[6,11] = 0
$data[6,12] = 0
$data[6,13] = 3.6466955528673
$data[7,0] = 0.514862142623798
$data[7,1] = 0.05544724859531414
$data[7,2] = 0.6799660935058967
$data[7,3] = 0.2793517597284136
$data[7,4] = 0
$data[7,5] = 0.8631039296232359
$data[7,6] = 0.9234828015470526
$data[7,7] = 0
$data[7,8] = 0.1470878135779543
$data[7,9] = 0.4874800762892733
$data[7,10] = 0
$data[7,11] = 0
$data[7,12] = 0
$data[7,13] = 3.601707505546585
$data[8,0] = 0.5889469993342686
$data[8,1] = 0.06500188873559409
$data[8,2] = 0.7021160887527174
$data[8,3] = 0.2900239945673064
$data[8,4] = 0
$data[8,5] = 0.8600054058292983
$data[8,6] = 0.9152067845807608
$data[8,7] = 0
$data[8,8] = 0.1540161720966324
$data[8,9] = 0.5659104058546518
$data[8,10] = 0
$data[8,11] = 0
$data[8,12] = 0
$data[8,13] = 3.578277549998489
$data[9,0] = 0.6227044247654874
$data[9,1] = 0.06933032602371725
$data[9,2] = 0.7125387684589271
$data[9,3] = 0.2950274491293001
$data[9,4] = 0
$data[9,5] = 0.8591719407251617
$data[9,6] = 0.9118634601834827
$data[9,7] = 0
$data[9,8] = 0.1572504753011543
$data[9,9] = 0.6015842650039644
$data[9,10] = 0
$data[9,11] = 0
$data[9,12] = 0
$data[9,13] = 3.56971178481416
$data[10,0] = 0.6354950468183631
$data[10,1] = 0.07096673344251769
$data[10,2] = 0.7165353120042823
$data[10,3] = 0.2969434800708157
$data[10,4] = 0
$data[10,5] = 0.8589393167912789
$data[10,6] = 0.9106580085614411
$data[10,7] = 0
$data[10,8] = 0.1584871028853883
$data[10,9] = 0.6150918834852064
$data[10,10] = 0
$data[10,11] = 0
$data[10,12] = 0
$data[10,13] = 3.566769334098723
$data[11,0] = 0.6327400373494072
$data[11,1] = 0.070614424303443
$data[11,2] = 0.7156723764607023
$data[11,3] = 0.2965298804077889
$data[11,4] = 0
$data[11,5] = 0.85898572302753
$data[11,6] = 0.910914929492975
$data[11,7] = 0
$data[11,8] = 0.15822024528903
$data[11,9] = 0.6121828442260551
$data[11,10] = 0
$data[11,11] = 0
$data[11,12] = 0
$data[11,13] = 3.567389641957107
$data[12,0] = 0.6237565717850941
$data[12,1] = 0.06946500841161196
$data[12,2] = 0.7128665711684619
$data[12,3] = 0.2951846548518233
$data[12,4] = 0
$data[12,5] = 0.859151138471276
$data[12,6] = 0.9117630726687906
$data[12,7] = 0
$data[12,8] = 0.1573519755880284
$data[12,9] = 0.6026955746371527
$data[12,10] = 0
$data[12,11] = 0
$data[12,12] = 0
$data[12,13] = 3.569463669037162
$data[13,0] = 0.6182548875621876
$data[13,1] = 0.06876060668112416
$data[13,2] = 0.7111544030667574
$data[13,3] = 0.2943634422073771
$data[13,4] = 0
$data[13,5] = 0.8592632726491587
$data[13,6] = 0.9122904753790664
$data[13,7] = 0
$data[13,8] = 0.1568216806407037
$data[13,9] = 0.5968841597767209
$data[13,10] = 0
$data[13,11] = 0
$data[13,12] = 0
$data[13,13] = 3.5707733069585
$data[14,0] = 0.5867419435537329
$data[14,1] = 0.0647186451820545
$data[14,2] = 0.7014419041567805
$data[14,3] = 0.2896999940895242
$data[14,4] = 0
$data[14,5] = 0.8600714798233042
$data[14,6] = 0.9154337585876675
$data[14,7] = 0
$data[14,8] = 0.1538064636664416
$data[14,9] = 0.5635788890835443
$data[14,10] = 0
$data[14,11] = 0
$data[14,12] = 0
$data[14,13] = 3.578879477774393
$data[15,0] = 0.56742364067415
$data[15,1] = 0.06223435609689432
$data[15,2] = 0.6955722639842179
$data[15,3] = 0.286877153520571
$data[15,4] = 0
$data[15,5] = 0.8607149430139174
$data[15,6] = 0.9174700005390406
$data[15,7] = 0
$data[15,8] = 0.1519778675987027
$data[15,9] = 0.5431455690624603
$data[15,10] = 0
$data[15,11] = 0
$data[15,12] = 0
$data[15,13] = 3.584388539464669
$data[16,0] = 0.5563175574378647
$data[16,1] = 0.06080376809852339
$data[16,2] = 0.6922288276014683
$data[16,3] = 0.2852675196007795
$data[16,4] = 0
$data[16,5] = 0.8611392592280822
$data[16,6] = 0.9186808653552276
$data[16,7] = 0
$data[16,8] = 0.1509338785556906
$data[16,9] = 0.5313924798582832
$data[16,10] = 0
$data[16,11] = 0
$data[16,12] = 0
$data[16,13] = 3.587754157219592
$data[17,0] = 0.5525581564163531
$data[17,1] = 0.06031910787868355
$data[17,2] = 0.6911024049067294
$data[17,3] = 0.2847249290286129
$data[17,4] = 0
$data[17,5] = 0.8612922319025529
$data[17,6] = 0.919097657604425
$data[17,7] = 0
$data[17,8] = 0.1505817366270463
$data[17,9] = 0.5274130397716874
$data[17,10] = 0
$data[17,11] = 0
$data[17,12] = 0
$data[17,13] = 3.588927514990701
$data[18,0] = 0.5694795635796197
$data[18,1] = 0.062498988519053
$data[18,2] = 0.696193721540368
$data[18,3] = 0.2871762021479043
$data[18,4] = 0
$data[18,5] = 0.8606408332821758
$data[18,6] = 0.9172491332820272
$data[18,7] = 0
$data[18,8] = 0.1521717205308875
$data[18,9] = 0.5453207753342042
$data[18,10] = 0
$data[18,11] = 0
$data[18,12] = 0
$data[18,13] = 3.583781704555463
$data[19,0] = 0.6263950369499014
$data[19,1] = 0.06980269308864706
$data[19,2] = 0.7136893563681781
$data[19,3] = 0.2955792015802103
$data[19,4] = 0
$data[19,5] = 0.8591002983702367
$data[19,6] = 0.9115123081072767
$data[19,7] = 0
$data[19,8] = 0.1576066854125173
$data[19,9] = 0.6054822560901698
$data[19,10] = 0
$data[19,11] = 0
$data[19,12] = 0
$data[19,13] = 3.56884629963389
$data[20,0] = 0.6636355418332869
$data[20,1] = 0.07456043775442822
$data[20,2] = 0.7254133996164569
$data[20,3] = 0.3011953650966177
$data[20,4] = 0
$data[20,5] = 0.858577276730486
$data[20,6] = 0.9081161355326373
$data[20,7] = 0
$data[20,8] = 0.161227917787059
$data[20,9] = 0.6447933970635802
$data[20,10] = 0
$data[20,11] = 0
$data[20,12] = 0
$data[20,13] = 3.560841020833038
$data[21,0] = 0.6437558554640361
$data[21,1] = 0.07202260033321295
$data[21,2] = 0.7191295986513637
$data[21,3] = 0.2981865502514225
$data[21,4] = 0
$data[21,5] = 0.8588121043861605
$data[21,6] = 0.9098964261682454
$data[21,7] = 0
$data[21,8] = 0.1592888705333309
$data[21,9] = 0.6238132354453114
$data[21,10] = 0
$data[21,11] = 0
$data[21,12] = 0
$data[21,13] = 3.564952827487957
$data[22,0] = 0.5685500799588112
$data[22,1] = 0.06237935547531492
$data[22,2] = 0.6959126637037514
$data[22,3] = 0.2870409609747568
$data[22,4] = 0
$data[22,5] = 0.8606741689093553
$data[22,6] = 0.917348862137672
$data[22,7] = 0
$data[22,8] = 0.1520840568972233
$data[22,9] = 0.5443373822899957
$data[22,10] = 0
$data[22,11] = 0
$data[22,12] = 0
$data[22,13] = 3.584055436674333
$data[23,0] = 0.4876330648266389
$data[23,1] = 0.05191323968800532
$data[23,2] = 0.672117184693235
$data[23,3] = 0.2755535730894252
$data[23,4] = 0
$data[23,5] = 0.8647547284020476
$data[23,6] = 0.9269043393222063
$data[23,7] = 0
$data[23,8] = 0.1446095844466342
$data[23,9] = 0.4585974783697395
$data[23,10] = 0
$data[23,11] = 0
$data[23,12] = 0
$data[23,13] = 3.61218924813025

$ws.Range("B2:O25").Value = $data
